$d = $word.ActiveDocument

# The new sub-bullet items to append after "Logique 2 : ..." (week 3 update).
$newItems = @(
    "Reflexion -> logique 2, compliqué et surchage de boucle pour comparé dans la base de versionning ce qui a changer etc.",
    "J'ai choisi la logique 1, apres avoir trouver une api java qui compare 2 objets, et renvoie les differences entre les 2, ce qui permet de savoir ce qui a reelement changer, et nous permettra d'ajouter directement dans la base de versionning que les champs qui ont été modifié.",
    "Test de l'api java Javers qui compare 2 objet et renvoie les champs qui diffèrent.",
    "Implémentation de Javers dans le projet"
)

# Start from the last paragraph in the document ("Logique 2 : ..."), which
# carries the pStyle "Paragraphedeliste" + numPr ilvl=1/numId=1 that the new
# paragraphs must reuse.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

foreach ($text in $newItems) {
    $r = $target.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
    $target.Range.InsertAfter($text)
}
